$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Activités")
$ws2 = $wb.Worksheets.Item("Données")

# --- Row 21: new journal entry ---
$ws1.Range("A21").Value = 44265
$ws1.Range("B21").Value = 0.33611111111111108
$ws1.Range("C21").Value = 0.41805555555555557
$ws1.Range("G21").Value = "Remise à jour du nouveau Github"
$ws1.Range("E21").Value = "Rédaction documentation"

# --- Row 22: new journal entry ---
$ws1.Range("A22").Value = 44265
$ws1.Range("B22").Value = 0.41805555555555557
$ws1.Range("C22").Value = 0.44097222222222227
$ws1.Range("G22").Value = "Avancement sur le jeu"
$ws1.Range("E22").Value = "Création"

# --- Données sheet: new lookup row matching the new activity ---
$ws2.Range("A10").Value = "Rédaction documentation"
$ws2.Range("B10").Value = "Documentation"

# --- Update selections / active view to match the final saved state ---
$ws2.Activate()
$ws2.Range("D14").Select()

$ws1.Activate()
$ws1.Range("G26").Select()
